$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New accelerometer samples added on may 9th (rows 22-31, columns x/y/z)
$data = @(
    @(22, -3.367526054382324, 0.915987193584442, -1.626443386077881),
    @(23, -1.118759155273438, 2.792432069778442, -4.963344097137451),
    @(24, -1.557756900787354, -0.5582034587860107, -0.2619988918304443),
    @(25, 1.715949058532715, -1.576748490333557, 5.096891403198242),
    @(26, -0.8243503570556641, 0.5943599939346313, 1.927432060241699),
    @(27, -0.2981023788452148, 1.024843096733093, 0.8517363667488098),
    @(28, 0.4371089935302734, 0.3337190449237823, -0.154114544391632),
    @(29, -0.5641984939575195, -0.3292053341865539, -0.326197862625122),
    @(30, 0.131052017211914, 0.5107872486114502, 0.0540084838867187),
    @(31, -0.1777238845825195, 0.4102384448051452, 0.1352127194404602)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
}
